# Apply the LinuxForHealth re-brand / version bump edit to the FHIR
# StructureDefinition workbook.
#
# Net effect (verified against the target diff):
#   Metadata!B2  URL        -> http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-id
#   Metadata!B3  Version    -> 8.0.0
#   Metadata!B8  Date       -> 2022-11-10T16:00:46+00:00
#   Metadata!B9  Publisher  -> LinuxForHealth Team
#   Elements!Q5  (Extension.url Fixed Value) -> same new URL as above
#   Elements!AI2 (Extension row, Constraint(s)) -> cleared
#
# (The underlying xlsx diff also reshuffles the shared-string table so that
# the "Extension.id" / "Extension.extension" rows' strings sit ahead of the
# ele-1/ext-1 constraint text, and the constraint text itself is only kept
# on the Extension.extension row (AI4) instead of being duplicated on the
# Extension row (AI2) -- both rows pointed at the identical shared string
# before the edit, so clearing AI2 is the only value-level change needed.)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-id"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-id"
$wsElem.Range("AI2").Value = ""
